# GMS Data Release 1
# Updates the Tiering_data dictionary sheet:
#  - rename field "rare_diseases_family_id" -> "referral_id"
#  - remove the "sample_id" field row entirely
#  - rename field "assembly" -> "genome_build"
#  - merge "full_brothers_affected"/"full_sisters_affected" rows into a
#    single "full_siblings_affected" row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename rare_diseases_family_id -> referral_id (row 3, column B)
$ws.Range("B3").Value = "referral_id"

# Remove the sample_id row (row 6) completely; remaining rows shift up.
$ws.Rows(6).Delete()

# Merge full_brothers_affected (row 24) / full_sisters_affected (row 25)
# into a single full_siblings_affected row, then delete the now-duplicate row.
$ws.Range("B24").Value = "full_siblings_affected"
$ws.Range("D24").Value = "Number of full siblings with same condition"
$ws.Rows(25).Delete()

# Rename assembly -> genome_build (row 8)
$ws.Range("B8").Value = "genome_build"
